# Add a new "runner_settings" column immediately before the existing "id"
# column on the TestAsset, AcceptanceTestAsset and TestEdgeData sheets.
#
# On each of these sheets the trailing columns are: ... id | name | description | tags
# After the edit they become:                       ... runner_settings | id | name | description | tags

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> column letter where "id" currently lives (i.e.
# the column a new blank column needs to be inserted in front of).
$targets = @{
    "TestAsset"           = "K"
    "AcceptanceTestAsset" = "U"
    "TestEdgeData"        = "K"
}

foreach ($sheetName in $targets.Keys) {
    $idColumn = $targets[$sheetName]
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a new blank column in front of the "id" column, shifting
    # "id", "name", "description" and "tags" one column to the right.
    $ws.Columns($idColumn).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftToRight)

    # Populate the header for the newly inserted column.
    $ws.Range($idColumn + "1").Value = "runner_settings"
}
